$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author added two new price records to the top of the data block that
# starts at row 24, pushing the previously-existing rows (24-75) down by two
# rows (to 26-77). Inserting whole rows preserves all the existing data,
# formatting (e.g. the date style on column D) and simply shifts it down.
$ws.Rows("24:25").Insert()

# New row 24: Región de La Araucanía entry on 44536 (2021-12-06)
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44536
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 100112022
$ws.Range("G24").Value = "Arveja Verde"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 18000
$ws.Range("L24").Value = 18000
$ws.Range("M24").Value = 18000
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Región de La Araucanía"
$ws.Range("P24").Value = 720
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"

# New row 25: Región del Maule entry also on 44536 (2021-12-06)
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44536
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = 100112022
$ws.Range("G25").Value = "Arveja Verde"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 15000
$ws.Range("N25").Value = "$/saco 25 kilos"
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 600
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"

# Make sure the date cells keep the workbook's date/time number format,
# matching the style used by the rest of column D (style index "2" ->
# numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
